# "Generate Report for Handback"
#
# The handback transform for the 36bc403d... file failed (file name
# mismatch between the handback package and the original handoff), so
# the localization-status report needs to reflect that failure instead
# of the previous "Ready for handoff" status, and record the error
# detail message for both the zh-cn and de-de targets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Status for the 36bc403d... row changes everywhere it is surfaced:
# the per-locale sheets and the roll-up Overview sheet.
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$zhcn.Range("C3").Value     = $newStatus
$dede.Range("C3").Value     = $newStatus

# Record the handback/handoff file-name-mismatch error detail for the
# 36bc403d... row on each locale sheet.
$zhcn.Range("P3").Value = "Handback file name: gvrmcmlm.hpp is different with handoff file name: 36bc403d-c79e-42e3-b5af-727a3793fbaa.0c6b7185ed40aac1a3973946a8e4c1bc251fd9c8.zh-cn."
$dede.Range("P3").Value = "Handback file name: gvrmcmlm.hpp is different with handoff file name: 36bc403d-c79e-42e3-b5af-727a3793fbaa.0c6b7185ed40aac1a3973946a8e4c1bc251fd9c8.de-de."

# Widen the "Error Detail" column (P, the 16th column) on both locale
# sheets so the new, much longer error messages are readable.
$zhcn.Columns.Item(16).ColumnWidth = 39.17
$dede.Columns.Item(16).ColumnWidth = 39.17
